# "changes in lead controller"
# Restructure the attendee sheet: insert Phone Number/Enrollment/Year/Branch
# columns after Email, move College to a new trailing column, and append a
# new attendee row (Het Patel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Phone Number"
$ws.Range("D1").Value = "Enrollment"
$ws.Range("E1").Value = "Year"
$ws.Range("F1").Value = "Branch"
$ws.Range("G1").Value = "College"

# --- Row 2: Kandarp Shah ---------------------------------------------------
$ws.Range("A2").Value = "Kandarp Shah"
$ws.Range("B2").Value = "shahkandarp24@gmail.com"
$ws.Range("C2").Value = "'7016763640"
$ws.Range("D2:F2").ClearContents()
$ws.Range("G2").Value = "GCET"

# --- Row 3: Pratham Shah ---------------------------------------------------
$ws.Range("A3").Value = "Pratham Shah"
$ws.Range("B3").Value = "prathamshah019@gmail.com"
$ws.Range("C3").Value = "'7405802474"
$ws.Range("D3:F3").ClearContents()
$ws.Range("G3").Value = "GCET"

# --- Row 4: Het Patel (new attendee) ---------------------------------------
$ws.Range("A4").Value = "Het Patel"
$ws.Range("B4").Value = "hetpatel5542@gmail.com"
$ws.Range("C4").Value = "'7698545581"
$ws.Range("D4:F4").ClearContents()
$ws.Range("G4").Value = "GCET"
